$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-09-04 12:43:22"

$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-09-04 12:43:17"

$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-09-04 12:43:22"

$wsOverview.Columns.Item(5).ColumnWidth = 16.333333333333334
$wsOverview.Columns.Item(6).ColumnWidth = 16.333333333333334
$wsZhCn.Columns.Item(3).ColumnWidth = 16.333333333333334
$wsDeDe.Columns.Item(3).ColumnWidth = 16.333333333333334
